$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 13158284
$ws.Range("I107").Value = 16667079
$ws.Range("J107").Value = 304
$ws.Range("K107").Value = 16667079
$ws.Range("L107").Value = 304
$ws.Range("M107").Value = -16665159
$ws.Range("N107").Value = -4144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 571.96295
$ws.Range("I2").Value = 422.91666
$ws.Range("K2").Value = 422.91666
$ws.Range("M2").Value = -309.91666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 280961.16
$ws.Range("I61").Value = 2497.8235
$ws.Range("J61").Value = 530112.5600000001
$ws.Range("K61").Value = 2497.8235
$ws.Range("L61").Value = 530112.5600000001
$ws.Range("M61").Value = -2285.8235
$ws.Range("N61").Value = -530536.5600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 906.1429000000001
$ws.Range("I97").Value = 687.1429000000001
$ws.Range("J97").Value = 1125.1428
$ws.Range("K97").Value = 687.1429000000001
$ws.Range("L97").Value = 1125.1428
$ws.Range("M97").Value = -191.1429000000001
$ws.Range("N97").Value = -2117.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2851326.8
$ws.Range("I102").Value = 3369268
$ws.Range("K102").Value = 3369268
$ws.Range("M102").Value = -3367646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2513.15
$ws.Range("I110").Value = 3121.3635
$ws.Range("K110").Value = 3121.3635
$ws.Range("M110").Value = -1076.3635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 571.96295
$ws.Range("I116").Value = 422.91666
$ws.Range("K116").Value = 422.91666
$ws.Range("M116").Value = 1871.08334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 755612.5
$ws.Range("I122").Value = 778482.6
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 2335447.8
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -2332997.8
$ws.Range("N122").Value = -7600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2780988.5
$ws.Range("I132").Value = 2079.7307
$ws.Range("J132").Value = 10006151
$ws.Range("K132").Value = 6239.1921
$ws.Range("L132").Value = 30018453
$ws.Range("M132").Value = -3709.1921
$ws.Range("N132").Value = -30023513

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 280961.16
$ws.Range("I136").Value = 2497.8235
$ws.Range("J136").Value = 530112.5600000001
$ws.Range("K136").Value = 7493.470499999999
$ws.Range("L136").Value = 1590337.68
$ws.Range("M136").Value = -4943.470499999999
$ws.Range("N136").Value = -1595437.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 571.96295
$ws.Range("I3").Value = 422.91666
$ws.Range("K3").Value = 422.91666
$ws.Range("M3").Value = -308.91666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12837.818
$ws.Range("I20").Value = 1272.1333
$ws.Range("J20").Value = 37621.43
$ws.Range("K20").Value = 1272.1333
$ws.Range("L20").Value = 37621.43
$ws.Range("M20").Value = -1025.1333
$ws.Range("N20").Value = -38115.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2050.3125
$ws.Range("I86").Value = 2000.5555
$ws.Range("J86").Value = 2114.2856
$ws.Range("K86").Value = 2000.5555
$ws.Range("L86").Value = 2114.2856
$ws.Range("M86").Value = -877.5554999999999
$ws.Range("N86").Value = -4360.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2050.3125
$ws.Range("I89").Value = 2000.5555
$ws.Range("J89").Value = 2114.2856
$ws.Range("K89").Value = 10002.7775
$ws.Range("L89").Value = 10571.428
$ws.Range("M89").Value = -4386.7775
$ws.Range("N89").Value = -21803.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 100001180
$ws.Range("I99").Value = 200000800
$ws.Range("J99").Value = 1556
$ws.Range("K99").Value = 200000800
$ws.Range("L99").Value = 1556
$ws.Range("M99").Value = -199999302
$ws.Range("N99").Value = -4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11278.409
$ws.Range("I105").Value = 16043.214
$ws.Range("J105").Value = 2940
$ws.Range("K105").Value = 16043.214
$ws.Range("L105").Value = 2940
$ws.Range("M105").Value = -14296.214
$ws.Range("N105").Value = -6434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1168.1578
$ws.Range("I107").Value = 990.5
$ws.Range("K107").Value = 990.5
$ws.Range("M107").Value = 929.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1760.4814
$ws.Range("I16").Value = 1435.7059
$ws.Range("J16").Value = 2312.6
$ws.Range("K16").Value = 1435.7059
$ws.Range("L16").Value = 2312.6
$ws.Range("M16").Value = -1148.7059
$ws.Range("N16").Value = -2886.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4300.4707
$ws.Range("I31").Value = 1909.1765
$ws.Range("J31").Value = 6691.7646
$ws.Range("K31").Value = 1909.1765
$ws.Range("L31").Value = 6691.7646
$ws.Range("M31").Value = -1614.1765
$ws.Range("N31").Value = -7281.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4300.4707
$ws.Range("I34").Value = 1909.1765
$ws.Range("J34").Value = 6691.7646
$ws.Range("K34").Value = 1909.1765
$ws.Range("L34").Value = 6691.7646
$ws.Range("M34").Value = -1707.1765
$ws.Range("N34").Value = -7095.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1834.3334
$ws.Range("I105").Value = 2084.9167
$ws.Range("J105").Value = 1333.1666
$ws.Range("K105").Value = 2084.9167
$ws.Range("L105").Value = 1333.1666
$ws.Range("M105").Value = -337.9167000000002
$ws.Range("N105").Value = -4827.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1135.65
$ws.Range("I107").Value = 1116.3636
$ws.Range("J107").Value = 1159.2222
$ws.Range("K107").Value = 1116.3636
$ws.Range("L107").Value = 1159.2222
$ws.Range("M107").Value = 803.6364000000001
$ws.Range("N107").Value = -4999.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1760.4814
$ws.Range("I113").Value = 1435.7059
$ws.Range("J113").Value = 2312.6
$ws.Range("K113").Value = 1435.7059
$ws.Range("L113").Value = 2312.6
$ws.Range("M113").Value = 734.2941000000001
$ws.Range("N113").Value = -6652.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3026.875
$ws.Range("I122").Value = 459
$ws.Range("J122").Value = 3772.3872
$ws.Range("K122").Value = 4131
$ws.Range("L122").Value = 33951.4848
$ws.Range("M122").Value = -1681
$ws.Range("N122").Value = -38851.4848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 58825110
$ws.Range("I113").Value = 83334376
$ws.Range("K113").Value = 83334376
$ws.Range("M113").Value = -83332206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3404.5356
$ws.Range("I132").Value = 4101.5
$ws.Range("J132").Value = 2707.5715
$ws.Range("K132").Value = 12304.5
$ws.Range("L132").Value = 8122.7145
$ws.Range("M132").Value = -9774.5
$ws.Range("N132").Value = -13182.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1593.6154
$ws.Range("I61").Value = 1583.3636
$ws.Range("J61").Value = 1650
$ws.Range("K61").Value = 1583.3636
$ws.Range("L61").Value = 1650
$ws.Range("M61").Value = -1381.3636
$ws.Range("N61").Value = -2054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2161
$ws.Range("I93").Value = 2453.2
$ws.Range("K93").Value = 2453.2
$ws.Range("M93").Value = -1205.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1105.875
$ws.Range("I100").Value = 978.1429000000001
$ws.Range("K100").Value = 978.1429000000001
$ws.Range("M100").Value = -437.1429000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1593.6154
$ws.Range("I113").Value = 1583.3636
$ws.Range("J113").Value = 1650
$ws.Range("K113").Value = 1583.3636
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 586.6364000000001
$ws.Range("N113").Value = -5990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 66333.336
$ws.Range("J127").Value = 66333.336
$ws.Range("L127").Value = 66333.336
$ws.Range("N127").Value = -76253.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3990.9465
$ws.Range("I136").Value = 2067.9546
$ws.Range("J136").Value = 11041.917
$ws.Range("K136").Value = 6203.8638
$ws.Range("L136").Value = 33125.751
$ws.Range("M136").Value = -3653.8638
$ws.Range("N136").Value = -38225.751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 84150.164
$ws.Range("I100").Value = 200560.4
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 401120.8
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -400579.8
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 500000640
$ws.Range("I107").Value = 500000640
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1500001920
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1500000000
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 902.88464
$ws.Range("I113").Value = 803
$ws.Range("J113").Value = 3400
$ws.Range("K113").Value = 2409
$ws.Range("L113").Value = 10200
$ws.Range("M113").Value = -239
$ws.Range("N113").Value = -14540

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1325.25
$ws.Range("I122").Value = 1000.3333
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 3000.9999
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -550.9998999999998
$ws.Range("N122").Value = -11800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2671.3125
$ws.Range("I132").Value = 2800
$ws.Range("J132").Value = 2612.818
$ws.Range("K132").Value = 8400
$ws.Range("L132").Value = 7838.454000000001
$ws.Range("M132").Value = -5870
$ws.Range("N132").Value = -12898.454
